$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.210687756538391
$ws.Range("B1").Value = 2.598259687423706
$ws.Range("C1").Value = 9.375737190246582
$ws.Range("D1").Value = 2.050379276275635
$ws.Range("E1").Value = 1.195765614509583
